$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Receptor-related recalculated TPM values)
$ws.Range("M2").Value = 0.08962966666666666
$ws.Range("N2").Value = 0.268889
$ws.Range("O2").Value = 0.4339761198462219
$ws.Range("P2").Value = 0.4339761198462219
$ws.Range("Q2").Value = 0.01420334438766667
$ws.Range("R2").Value = 0.127830099489
$ws.Range("S2").Value = 0.4339761198462219
$ws.Range("T2").Value = 0.4339761198462219

# Row 3 updates (Receptor-related recalculated TPM values)
$ws.Range("M3").Value = 0.1169016666666667
$ws.Range("N3").Value = 0.350705
$ws.Range("O3").Value = 0.5660238801537781
$ws.Range("P3").Value = 0.5660238801537781
$ws.Range("Q3").Value = 0.01852505641166667
$ws.Range("R3").Value = 0.166725507705
$ws.Range("S3").Value = 0.5660238801537781
$ws.Range("T3").Value = 0.5660238801537781
